# NC 2016 statewide alpha and beta calculation.xlsx
# "fixed NC alphas and betas in input file"
#
# Summary of the edit applied here:
#  1. Sheet "Computing alpha and beta":
#     - "Raw data based on:" -> "Raw data from:"
#     - "column sum" -> "row sum", "divide by column sum" -> "divide by row sum"
#     - Narrow column A
#     - Add two explanatory notes near the bottom (A57, B58)
#     - Move the view/selection away from this sheet (no longer the tab shown)
#  2. Sheet "For input to Pyro":
#     - Fix the mislabeled/garbled alpha rows: rotate E2:E4 back into the
#       correct order and drop the stray number formatting on E2
#     - Make this sheet the selected/active one, with E3 selected

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Computing alpha and beta")
$ws2 = $wb.Worksheets.Item("For input to Pyro")

# ---------------------------------------------------------------------------
# 1. Sheet 1 ("Computing alpha and beta") text / label fixes
# ---------------------------------------------------------------------------

# "Raw data based on:" -> "Raw data from:"
$ws1.Range("A1").Value = "Raw data from:"

# "column sum" -> "row sum"
$ws1.Range("E42").Value = "row sum"

# "divide by column sum" -> "divide by row sum"
$ws1.Range("C47").Value = "divide by row sum"

# New explanatory notes added near the bottom of the sheet
$ws1.Range("A57").Value = "In the input tab, the candidates are: 0 (other), 1 (Dem), 2 (Rep)"
$ws1.Range("B58").Value = "the races are: 0 (white), 1 (black), 2(other)"

# Column A is narrower now
$ws1.Columns.Item(1).ColumnWidth = 11.75

# ---------------------------------------------------------------------------
# 2. Sheet 2 ("For input to Pyro") data fix: rotate the alpha (E2:E4) values
#    back into their correct order, and drop the stray decimal formatting
#    that had been applied to E2.
# ---------------------------------------------------------------------------

$ws2.Range("E2").Style = "Normal"
$ws2.Range("E2").Value = 0.28000000000000003
$ws2.Range("E3").Value = 0.2
$ws2.Range("E4").Value = -0.48

# ---------------------------------------------------------------------------
# 3. View-state changes: sheet 1 is no longer the active/selected tab;
#    sheet 2 is now active, scrolled/selected at E3. Sheet 1's own
#    scroll position/selection also moved (topLeftCell A22->A21,
#    selection H49->B61).
# ---------------------------------------------------------------------------

($ws1.Activate()) | Out-Null
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
($ws1.Range("B61").Select()) | Out-Null

($ws2.Activate()) | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
($ws2.Range("E3").Select()) | Out-Null

# Window geometry tweak recorded alongside the tab switch
$excel.ActiveWindow.Left = 800
$excel.ActiveWindow.Top = 460
$excel.ActiveWindow.Width = 11140
$excel.ActiveWindow.Height = 16660
